$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    # Force text interpretation so numeric-looking strings (e.g. "1.004")
    # are not auto-converted to numbers by Excel type inference, then
    # drop the temporary number-format style so the cell keeps its original
    # (unstyled) formatting -- matches how these data cells were authored.
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextCell 'D2' '29.077.75'
Set-TextCell 'E2' '  -1.68%  '

Set-TextCell 'D3' '1.966.16'
Set-TextCell 'E3' '  -2.87%  '

Set-TextCell 'D4' '1.004'
Set-TextCell 'E4' '  -1.46%  '

Set-TextCell 'D5' '327.88'
Set-TextCell 'E5' '  -1.28%  '

Set-TextCell 'D6' '1.005'
Set-TextCell 'E6' '  -0.89%  '

Set-TextCell 'D7' '0.4984'
Set-TextCell 'E7' '  +0.85%  '

Set-TextCell 'D8' '0.4217'
Set-TextCell 'E8' '  +1.25%  '

Set-TextCell 'D9' '52.77'
Set-TextCell 'E9' '  -1.61%  '

Set-TextCell 'D10' '0.09286'
Set-TextCell 'E10' '  +5.55%  '

Set-TextCell 'D11' '1.099'
Set-TextCell 'E11' '  -1.66%  '

Set-TextCell 'D12' '22.90'
Set-TextCell 'E12' '  +1.02%  '

Set-TextCell 'D13' '1.967.31'
Set-TextCell 'E13' '  -5.27%  '

Set-TextCell 'B14' 'Chainlink'
Set-TextCell 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D14' '7.882'
Set-TextCell 'E14' '  -3.10%  '

Set-TextCell 'B15' 'Polkadot'
Set-TextCell 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D15' '6.466'
Set-TextCell 'E15' '  -0.38%  '

Set-TextCell 'D16' '1.006'
Set-TextCell 'E16' '  -1.69%  '

Set-TextCell 'B17' 'Litecoin'
Set-TextCell 'C17' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D17' '91.73'
Set-TextCell 'E17' '  -5.47%  '

Set-TextCell 'B18' 'ShibaInu'
Set-TextCell 'C18' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D18' '0.00001102'
Set-TextCell 'E18' '  -0.04%  '

Set-TextCell 'D19' '0.06713'
Set-TextCell 'E19' '  +1.19%  '

Set-TextCell 'D20' '19.31'
Set-TextCell 'E20' '  +0.02%  '

Set-TextCell 'E21' '  -0.26%  '

Set-TextCell 'D22' '5.958'
Set-TextCell 'E22' '  -0.10%  '

Set-TextCell 'D23' '29.095.52'
Set-TextCell 'E23' '  -1.98%  '

Set-TextCell 'D24' '12.04'
Set-TextCell 'E24' '  +1.98%  '

Set-TextCell 'D25' '2.263'
Set-TextCell 'E25' '  -2.38%  '

Set-TextCell 'D26' '2.203.84'
Set-TextCell 'E26' '  -3.47%  '

Set-TextCell 'D27' '20.63'
Set-TextCell 'E27' '  +0.40%  '

Set-TextCell 'D28' '155.78'
Set-TextCell 'E28' '  -1.76%  '

Set-TextCell 'D29' '6.332'
Set-TextCell 'E29' '  -1.88%  '

Set-TextCell 'D30' '2.256'
Set-TextCell 'E30' '  -2.22%  '

Set-TextCell 'D31' '126.57'
Set-TextCell 'E31' '  -0.67%  '

Set-TextCell 'D32' '1.048'
Set-TextCell 'E32' '  +0.88%  '

Set-TextCell 'D33' '0.09842'
Set-TextCell 'E33' '  +0.73%  '

Set-TextCell 'D34' '1.528'
Set-TextCell 'E34' '  -0.23%  '

Set-TextCell 'D35' '5.829'
Set-TextCell 'E35' '  +0.58%  '

Set-TextCell 'D36' '3.678'
Set-TextCell 'E36' '  -4.98%  '

Set-TextCell 'D37' '0.02437'
Set-TextCell 'E37' '  +0.49%  '

Set-TextCell 'D38' '1.315'
Set-TextCell 'E38' '  +1.16%  '

Set-TextCell 'B39' 'FraxShare'
Set-TextCell 'C39' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D39' '9.044'
Set-TextCell 'E39' '  -6.58%  '

Set-TextCell 'B40' 'Hedera'
Set-TextCell 'C40' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D40' '0.06371'
Set-TextCell 'E40' '  +1.10%  '

Set-TextCell 'D41' '0.6463'
Set-TextCell 'E41' '  +0.30%  '

Set-TextCell 'D42' '11.43'
Set-TextCell 'E42' '  -2.67%  '

Set-TextCell 'D43' '0.1985'
Set-TextCell 'E43' '  -2.92%  '

Set-TextCell 'E44' '  -0.67%  '

Set-TextCell 'B45' 'Decentraland'
Set-TextCell 'C45' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D45' '0.6231'
Set-TextCell 'E45' '  -0.88%  '

Set-TextCell 'B46' 'WEMIXTOKEN'
Set-TextCell 'C46' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D46' '1.347'
Set-TextCell 'E46' '  +5.29%  '

Set-TextCell 'B47' 'NEARProtocol'
Set-TextCell 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D47' '2.208'
Set-TextCell 'E47' '  +1.00%  '

Set-TextCell 'B48' 'EnergySwap'
Set-TextCell 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D48' '13.39'
Set-TextCell 'E48' '  -1.59%  '

Set-TextCell 'D49' '3.469'
Set-TextCell 'E49' '  -3.23%  '

Set-TextCell 'D50' '0.00000000330'
Set-TextCell 'E50' '  -0.80%  '

Set-TextCell 'D51' '0.06980'
Set-TextCell 'E51' '  -1.41%  '

Write-Host "cryptos list updated"
